$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.102.29'
$ws.Range("E2").Value = '  +2.62%  '

$ws.Range("D3").Value = '2.351.16'

$ws.Range("E4").Value = '  +0.00%  '

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '545.09'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  +6.12%  '

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '134.91'
$cell.Style = $origStyle
$ws.Range("E6").Value = '  +2.42%  '

$ws.Range("E7").Value = '  +0.16%  '

$ws.Range("E8").Value = '  +0.81%  '

$ws.Range("D9").Value = '2.349.26'
$ws.Range("E9").Value = '  +1.10%  '

$ws.Range("E10").Value = '  +1.49%  '

$ws.Range("E11").Value = '  +1.19%  '

$ws.Range("E12").Value = '  +3.21%  '

$ws.Range("E13").Value = '  +6.68%  '

$ws.Range("D14").Value = '2.767.95'
$ws.Range("E14").Value = '  +1.41%  '

$ws.Range("E15").Value = '  +0.17%  '

$ws.Range("D16").Value = '58.075.60'
$ws.Range("E16").Value = '  +2.63%  '

$ws.Range("E17").Value = '  +0.69%  '

$ws.Range("D18").Value = '2.347.72'
$ws.Range("E18").Value = '  +1.39%  '

$ws.Range("E19").Value = '  +2.63%  '

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '335.11'
$cell.Style = $origStyle
$ws.Range("E20").Value = '  +2.13%  '

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.22'
$cell.Style = $origStyle
$ws.Range("E21").Value = '  +1.67%  '

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.70'
$cell.Style = $origStyle
$ws.Range("E22").Value = '  -0.55%  '

$ws.Range("E23").Value = '  +0.23%  '

$ws.Range("E24").Value = '  +0.86%  '

$ws.Range("E25").Value = '  +2.84%  '

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = $origStyle
$ws.Range("E26").Value = '  +0.27%  '

$ws.Range("E27").Value = '  -1.43%  '

$ws.Range("E28").Value = '  +7.93%  '

$ws.Range("E29").Value = '  +5.30%  '

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '170.51'
$cell.Style = $origStyle
$ws.Range("E30").Value = '  +1.56%  '

$ws.Range("D31").Value = '0.0₃0732'
$ws.Range("E31").Value = '  +1.94%  '

$ws.Range("E32").Value = '  +0.82%  '

$ws.Range("E33").Value = '  +17.69%  '

$ws.Range("E34").Value = '  +1.11%  '

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = $origStyle
$ws.Range("E36").Value = '  +0.32%  '

$ws.Range("E37").Value = '  +6.42%  '

$ws.Range("E38").Value = '  +1.67%  '

$ws.Range("E39").Value = '  +4.70%  '

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '39.38'
$cell.Style = $origStyle
$ws.Range("E40").Value = '  +2.01%  '

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '148.54'
$cell.Style = $origStyle
$ws.Range("E41").Value = '  -0.33%  '

$ws.Range("E42").Value = '  +1.61%  '

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '286.97'
$cell.Style = $origStyle
$ws.Range("E43").Value = '  +3.98%  '

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.60'
$cell.Style = $origStyle
$ws.Range("E44").Value = '  +1.08%  '

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '19.29'
$cell.Style = $origStyle
$ws.Range("E45").Value = '  +5.93%  '

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0928'
$cell.Style = $origStyle
$ws.Range("E46").Value = '  +0.05%  '

$ws.Range("E47").Value = '  +2.12%  '

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.564'
$cell.Style = $origStyle
$ws.Range("E48").Value = '  +1.83%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0217'
$cell.Style = $origStyle
$ws.Range("E49").Value = '  +1.44%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '17.59'
$cell.Style = $origStyle
$ws.Range("E50").Value = '  +2.98%  '

$ws.Range("B51").Value = 'Polygon'
$ws.Range("C51").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.382'
$cell.Style = $origStyle
$ws.Range("E51").Value = '  +1.02%  '
